# Je factuur wordt gegenereerd op basis van de bestelling
# Fill the invoice header (factuurnummer / factuurdatum / leverdatum) and
# append the ordered product lines + subtotal/BTW/totaal block, the way the
# webshop backend does when it generates an invoice from an order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Factuurnummer / Factuurdatum / Leverdatum (row 21) -------------------
# C21 used to hold the placeholder text "Het werkt" -> becomes the real
# (numeric) invoice number.
$ws.Range("C21").Value = 6

# F21 / I21 get a real date value plus a yyyy-mm-dd date format.
$ws.Range("F21").NumberFormat = "yyyy-mm-dd"
$ws.Range("F21").Value = 44733   # 2022-06-21 Factuurdatum

$ws.Range("I21").NumberFormat = "yyyy-mm-dd"
$ws.Range("I21").Value = 44738   # 2022-06-26 Leverdatum

# --- Ordered product lines (rows 24-26) ------------------------------------
$ws.Range("C24:G24").Merge()
$ws.Range("C24").Value = "Snackpan XL"
$ws.Range("H24").Value = "'1"
$ws.Range("I24").Value = "'35.00"
$ws.Range("J24").Value = "'9%"
$ws.Range("K24").Value = "'35.00"

$ws.Range("C25:G25").Merge()
$ws.Range("C25").Value = "Pan"
$ws.Range("H25").Value = "'2"
$ws.Range("I25").Value = "'15.00"
$ws.Range("J25").Value = "'9%"
$ws.Range("K25").Value = "'30.00"

$ws.Range("C26:G26").Merge()
$ws.Range("C26").Value = "Aardappel pan"
$ws.Range("H26").Value = "'1"
$ws.Range("I26").Value = "'15.00"
$ws.Range("J26").Value = "'9%"
$ws.Range("K26").Value = "'15.00"

# --- Subtotaal / BTW / Totaal (rows 27-29) ---------------------------------
$ws.Range("H27:I27").Merge()
$ws.Range("H27").Value = "Subtotaal"
$ws.Range("J27").Formula = "=J29-J28"

$ws.Range("H28:I28").Merge()
$ws.Range("H28").Value = "BTW"
$ws.Range("J28").Formula = "=J29*0.09"

$ws.Range("H29:I29").Merge()
$ws.Range("H29").Value = "Totaal"
$ws.Range("J29").Value = 80

Write-Output "invoice filled in from order"
